$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5885
$ws.Range("C3").Value = 11200
$ws.Range("D3").Value = 18100
$ws.Range("E3").Value = 28400
$ws.Range("F3").Value = 30500
$ws.Range("G3").Value = 25200
$ws.Range("B4").Value = 771.751936
$ws.Range("C4").Value = 1466.957824
$ws.Range("D4").Value = 2375.02464
$ws.Range("E4").Value = 3728.736256
$ws.Range("F4").Value = 3999.268864
$ws.Range("G4").Value = 3304.062976
$ws.Range("B5").Value = 168.7
$ws.Range("C5").Value = 174.75
$ws.Range("D5").Value = 218.15
$ws.Range("E5").Value = 255.3
$ws.Range("F5").Value = 503.94
$ws.Range("G5").Value = 1250.33
$ws.Range("B6").Value = 289
$ws.Range("C6").Value = 302
$ws.Range("D6").Value = 367
$ws.Range("E6").Value = 494
$ws.Range("F6").Value = 971
$ws.Range("G6").Value = 2802
$ws.Range("B7").Value = 302
$ws.Range("C7").Value = 338
$ws.Range("D7").Value = 416
$ws.Range("E7").Value = 578
$ws.Range("F7").Value = 1139
$ws.Range("G7").Value = 3621
$ws.Range("B12").Value = 27800
$ws.Range("C12").Value = 55400
$ws.Range("D12").Value = 107000
$ws.Range("E12").Value = 201000
$ws.Range("G12").Value = 296000
$ws.Range("B13").Value = 114.294784
$ws.Range("C13").Value = 226.492416
$ws.Range("D13").Value = 437.256192
$ws.Range("E13").Value = 823.13216
$ws.Range("F13").Value = 1251.999744
$ws.Range("G13").Value = 1212.153856
$ws.Range("B14").Value = 35.6589
$ws.Range("C14").Value = 35.71173
$ws.Range("D14").Value = 36.82257
$ws.Range("E14").Value = 38.90130000000001
$ws.Range("F14").Value = 51.47196
$ws.Range("G14").Value = 106.83795
$ws.Range("B15").Value = 62.72
$ws.Range("C15").Value = 64.256
$ws.Range("D15").Value = 69.12
$ws.Range("E15").Value = 78.336
$ws.Range("G15").Value = 342.016
$ws.Range("B16").Value = 91.648
$ws.Range("C16").Value = 92.672
$ws.Range("D16").Value = 95.744
$ws.Range("E16").Value = 102.912
$ws.Range("F16").Value = 136.192
$ws.Range("G16").Value = 505.856
$ws.Range("B21").Value = 5595
$ws.Range("C21").Value = 8865
$ws.Range("D21").Value = 13800
$ws.Range("E21").Value = 13700
$ws.Range("F21").Value = 15800
$ws.Range("G21").Value = 15900
$ws.Range("B22").Value = 732.954624
$ws.Range("C22").Value = 1161.822208
$ws.Range("D22").Value = 1814.03648
$ws.Range("E22").Value = 1793.06496
$ws.Range("F22").Value = 2075.131904
$ws.Range("G22").Value = 2089.811968
$ws.Range("B23").Value = 141.25
$ws.Range("C23").Value = 160.22
$ws.Range("D23").Value = 171.41
$ws.Range("E23").Value = 361.72
$ws.Range("F23").Value = 556.85
$ws.Range("G23").Value = 1197.49
$ws.Range("B24").Value = 163
$ws.Range("C24").Value = 200
$ws.Range("D24").Value = 221
$ws.Range("E24").Value = 482
$ws.Range("F24").Value = 873
$ws.Range("G24").Value = 1598
$ws.Range("B25").Value = 217
$ws.Range("C25").Value = 229
$ws.Range("D25").Value = 243
$ws.Range("E25").Value = 545
$ws.Range("F25").Value = 1029
$ws.Range("G25").Value = 1827
$ws.Range("B30").Value = 51600
$ws.Range("C30").Value = 101000
$ws.Range("D30").Value = 188000
$ws.Range("E30").Value = 265000
$ws.Range("F30").Value = 283000
$ws.Range("G30").Value = 317000
$ws.Range("B31").Value = 211.812352
$ws.Range("C31").Value = 415.236096
$ws.Range("D31").Value = 767.557632
$ws.Range("E31").Value = 1083.179008
$ws.Range("F31").Value = 1159.725056
$ws.Range("G31").Value = 1298.137088
$ws.Range("B32").Value = 17.33986
$ws.Range("C32").Value = 17.34603
$ws.Range("D32").Value = 17.7024
$ws.Range("E32").Value = 23.05
$ws.Range("F32").Value = 36.67
$ws.Range("G32").Value = 63.25
$ws.Range("B33").Value = 20.352
$ws.Range("C33").Value = 22.144
$ws.Range("D33").Value = 22.912
$ws.Range("E33").Value = 31.616
$ws.Range("F33").Value = 70
$ws.Range("G33").Value = 106
$ws.Range("B34").Value = 22.4
$ws.Range("C34").Value = 24.192
$ws.Range("D34").Value = 25.472
$ws.Range("E34").Value = 35.584
$ws.Range("F34").Value = 99
$ws.Range("G34").Value = 161
$ws.Range("B39").Value = 11900
$ws.Range("C39").Value = 15200
$ws.Range("D39").Value = 22000
$ws.Range("E39").Value = 21300
$ws.Range("F39").Value = 19700
$ws.Range("B40").Value = 1560.281088
$ws.Range("C40").Value = 1988.100096
$ws.Range("D40").Value = 2886.729728
$ws.Range("E40").Value = 2789.21216
$ws.Range("F40").Value = 2584.73984
$ws.Range("G40").Value = 2058.354688
$ws.Range("B41").Value = 83.53
$ws.Range("C41").Value = 115.77
$ws.Range("D41").Value = 175.91
$ws.Range("E41").Value = 353.93
$ws.Range("F41").Value = 781.41
$ws.Range("G41").Value = 2008.25
$ws.Range("B42").Value = 190
$ws.Range("C42").Value = 412
$ws.Range("D42").Value = 453
$ws.Range("E42").Value = 1106
$ws.Range("F42").Value = 2147
$ws.Range("G42").Value = 5407
$ws.Range("B43").Value = 202
$ws.Range("C43").Value = 519
$ws.Range("D43").Value = 553
$ws.Range("E43").Value = 1434
$ws.Range("F43").Value = 2900
$ws.Range("B48").Value = 318000
$ws.Range("C48").Value = 529000
$ws.Range("D48").Value = 679000
$ws.Range("E48").Value = 683000
$ws.Range("F48").Value = 646000
$ws.Range("G48").Value = 537000
$ws.Range("B49").Value = 1303.379968
$ws.Range("C49").Value = 2165.30944
$ws.Range("D49").Value = 2781.872128
$ws.Range("E49").Value = 2796.552192
$ws.Range("F49").Value = 2644.508672
$ws.Range("G49").Value = 2199.912448
$ws.Range("B50").Value = 2.92857
$ws.Range("C50").Value = 3.44758
$ws.Range("D50").Value = 5.54321
$ws.Range("E50").Value = 10.94433
$ws.Range("F50").Value = 23.40143
$ws.Range("G50").Value = 58.51288
$ws.Range("B51").Value = 1.688
$ws.Range("C51").Value = 2.096
$ws.Range("D51").Value = 4.448
$ws.Range("E51").Value = 6.176
$ws.Range("F51").Value = 6.496
$ws.Range("B52").Value = 80.384
$ws.Range("C52").Value = 103.936
$ws.Range("D52").Value = 166.912
$ws.Range("E52").Value = 403.456
$ws.Range("F52").Value = 1089.536
$ws.Range("B57").Value = 5417
$ws.Range("C57").Value = 8641
$ws.Range("D57").Value = 12700
$ws.Range("E57").Value = 14900
$ws.Range("F57").Value = 13600
$ws.Range("G57").Value = 14400
$ws.Range("B58").Value = 709.885952
$ws.Range("C58").Value = 1132.46208
$ws.Range("D58").Value = 1661.99296
$ws.Range("E58").Value = 1952.448512
$ws.Range("F58").Value = 1787.82208
$ws.Range("G58").Value = 1892.67968
$ws.Range("B59").Value = 134.6
$ws.Range("C59").Value = 139.2
$ws.Range("D59").Value = 145.88
$ws.Range("E59").Value = 197.17
$ws.Range("F59").Value = 561.77
$ws.Range("G59").Value = 1353.49
$ws.Range("B60").Value = 161
$ws.Range("C60").Value = 167
$ws.Range("D60").Value = 180
$ws.Range("E60").Value = 269
$ws.Range("F60").Value = 799
$ws.Range("B61").Value = 165
$ws.Range("C61").Value = 176
$ws.Range("D61").Value = 188
$ws.Range("E61").Value = 289
$ws.Range("F61").Value = 914
$ws.Range("G61").Value = 1860
$ws.Range("B66").Value = 54400
$ws.Range("C66").Value = 114000
$ws.Range("D66").Value = 175000
$ws.Range("E66").Value = 239000
$ws.Range("F66").Value = 362000
$ws.Range("G66").Value = 453000
$ws.Range("B67").Value = 223.346688
$ws.Range("C67").Value = 465.567744
$ws.Range("D67").Value = 715.128832
$ws.Range("E67").Value = 980.41856
$ws.Range("F67").Value = 1483.73504
$ws.Range("G67").Value = 1853.882368
$ws.Range("B68").Value = 16.51381
$ws.Range("C68").Value = 14.37554
$ws.Range("D68").Value = 17.8526
$ws.Range("E68").Value = 22.45
$ws.Range("F68").Value = 29.7
$ws.Range("G68").Value = 55
$ws.Range("B69").Value = 19.328
$ws.Range("C69").Value = 19.328
$ws.Range("D69").Value = 22.656
$ws.Range("E69").Value = 30.848
$ws.Range("F69").Value = 55
$ws.Range("G69").Value = 69
$ws.Range("B70").Value = 21.12
$ws.Range("C70").Value = 21.376
$ws.Range("D70").Value = 25.216
$ws.Range("E70").Value = 34.56
$ws.Range("F70").Value = 81
$ws.Range("G70").Value = 98
